$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Casos totales / activos / Recuperados / Muertes" figures for a handful of provinces.
$updates = @{
    "Valencia/Valencia"  = @(627, 7, 608, 12)
    "Alacant/Alicante"   = @(372, 11, 350, 11)
    "Caceres"            = @(206, 2, 194, 10)
    "Badajoz"            = @(91, 5, 86, 0)
    "Castello/Castellon" = @(104, 1, 102, 1)
}

# Data rows run from row 4 to row 60 (column A = province name).
$lastRow = 60
for ($r = 4; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($updates.ContainsKey($name)) {
        $vals = $updates[$name]
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
        $ws.Cells.Item($r, 5).Value = $vals[3]
    }
}

# Re-sort the table by "Casos totales" (column B) descending, as the source does after every update.
$dataRange = $ws.Range("A4:E$lastRow")
$sortKey = $ws.Range("B4:B$lastRow")
$dataRange.Sort($sortKey, 2)

# Bump the "last updated" timestamp shown in the title row.
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 11:46"
